# Update "F2E4" (想去人数 / want-to-go count) values across sheets, as published
# by the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 250
    3  = 2510
    7  = 1319
    8  = 1681
    9  = 179
    11 = 2342
    12 = 485
    13 = 152
    16 = 104
    17 = 94
    18 = 8529
    20 = 6595
    21 = 10708
    22 = 124
    23 = 187
    24 = 204
    25 = 292
    26 = 521
    27 = 192
    28 = 168
    29 = 2040
    30 = 52
    31 = 15
    32 = 4439
    33 = 372
    34 = 426
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    8  = 1169
    9  = 8
    18 = 16
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Range("F$row").Value = $sheet2Updates[$row]
}

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    4  = 250
    6  = 2510
    10 = 1319
    12 = 1681
    14 = 179
    15 = 2342
    17 = 485
    18 = 152
    22 = 104
    23 = 94
    24 = 8529
    26 = 6595
    27 = 10708
    28 = 8
    29 = 124
    30 = 187
    31 = 204
    32 = 292
    33 = 521
    37 = 192
    38 = 168
    39 = 15
    40 = 4439
    45 = 16
    47 = 426
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
